# Fix a typo in the title-slide credits box: the professor's first name
# was misspelled "Liz" and should read "Luz"
# (Profesora: Dra. Liz Roxana de Leon Lomeli -> ... Luz ...).
#
# Slide 1 (SlideID 261) / Shape 1 (Id 2, "CuadroTexto 1") holds the text.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

$target = "Liz"
$replacement = "Luz"

for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text -like "*$target*") {
        $run = $para.Runs(1, 1)
        $run.Text = $run.Text.Replace($target, $replacement)
    }
}
